# Add two more 3D-printed parts for the "Chest Reopening Phantom Concept 3"
# job to the "June 2018" print-request log sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("June 2018")

# Row 29 - Chest Reopening Concept 3 Midstops
$ws.Cells.Item(29, 1).Value = "28-06-2018"
$ws.Cells.Item(29, 2).Value = "28-06-2018"
$ws.Cells.Item(29, 3).Value = "Chest Reopening Concept 3 Midstops"
$ws.Cells.Item(29, 4).Value = 4
$ws.Cells.Item(29, 5).Value = "Polylite + PLA"
$ws.Cells.Item(29, 6).Value = 3
$ws.Cells.Item(29, 7).Value = 20
$ws.Cells.Item(29, 8).Value = 0.2
$ws.Cells.Item(29, 9).Value = "NA"

# Row 30 - Chest Reopening Concept 3 Left + Right Sternum
$ws.Cells.Item(30, 1).Value = "28-06-2018"
$ws.Cells.Item(30, 2).Value = "28-06-2018"
$ws.Cells.Item(30, 3).Value = "Chest Reopening Concept 3 Left + Right Sternum"
$ws.Cells.Item(30, 4).Value = 1
$ws.Cells.Item(30, 5).Value = "Polyflex"
$ws.Cells.Item(30, 6).Value = 2
$ws.Cells.Item(30, 7).Value = 20
$ws.Cells.Item(30, 8).Value = 0.4
$ws.Cells.Item(30, 9).Value = "NA"

# Column C ("Part") is a best-fit column - widen it now that it holds a
# longer string ("Chest Reopening Concept 3 Left + Right Sternum").
$ws.Columns.Item(3).ColumnWidth = 38.5

# Move/extend the active selection like the author's last recorded action.
$ws.Range("H31").Select()
